$d = $word.ActiveDocument

# --- 1) Remove the empty paragraph that directly precedes the
#        "Vor der Einzahlung abzutrennen" paragraph, and give that
#        paragraph a bottom border (single, 0.75pt, 1pt space, auto color).
$target = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text -like "*Vor der Einzahlung abzutrennen*") {
        $target = $i
        break
    }
}

if ($target -ne $null) {
    $prev = $d.Paragraphs.Item($target - 1)
    if ($prev.Range.Text.Trim().Length -eq 0) {
        $prev.Range.Delete()
        $target = $target - 1
    }

    $p = $d.Paragraphs.Item($target)
    $p.Borders.Item(-3).LineStyle = 1
    $p.Borders.Item(-3).LineWidth = 3
    $p.Borders.Item(-3).Color = -16777216
    $p.Borders.DistanceFromBottom = 1
}

# --- 2) Center the table and its rows.
$tbl = $d.Tables.Item(1)
$tbl.Alignment = 1
for ($i = 1; $i -le $tbl.Rows.Count; $i++) {
    $tbl.Rows.Item($i).Alignment = 1
}
